$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.829.23'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +3.17%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.881.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +3.34%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.29%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''327.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.50%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.28%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.4670'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.81%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.3942'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +2.57%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.07930'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +1.32%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.9773'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +2.01%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D12').Value = '''1.900.97'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +3.16%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''5.755'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +2.06%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''7.014'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +2.47%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''0.07000'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +2.03%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''88.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +2.57%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''1.007'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = '''0.00001012'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +1.78%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = '''  +1.95%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''1.005'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.33%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''28.827.57'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +2.99%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''5.366'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +1.01%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  +1.84%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.121'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +1.06%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''2.090.54'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  +1.19%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''153.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +1.38%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  +1.40%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''5.761'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +0.99%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''2.010'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +2.43%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''120.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +3.23%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.09401'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +1.58%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''0.9416'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +0.26%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''5.322'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +0.99%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  +3.32%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''3.353'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -1.97%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''0.05918'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -0.42%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.02121'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -1.19%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''1.150'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +0.40%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = '''  +4.93%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''0.5719'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +2.50%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''10.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +1.23%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.1794'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +1.57%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.07255'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +3.67%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''11.87'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +2.59%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''0.5345'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +1.88%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''1.154'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -6.91%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''2.130'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -4.21%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '''  +1.69%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''114.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +1.89%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''2.372'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  +3.04%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''1.007'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +0.55%  '
$ws.Range('E51').Style = 'Normal'
